# Add a new "Apropriação" row (date + hours) below the existing one.
#
# Existing sheet:
#   A1 "Data" | B1 "Quantidade de horas"
#   A2 01/10/13 (date)    | B2 02:00:00 (time)
#
# New row to append:
#   A3 02/10/13 (date, serial 41549) | B3 02:30:00 (time, fraction of day)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date cell - same display format used by the cell above (DD/MM/YY)
$ws.Range("A3").Value2 = 41549
$ws.Range("A3").NumberFormat = "DD/MM/YY"

# Time-of-day cell (2:30 -> 2.5 / 24 of a day) - same format as the cell above (HH:MM:SS)
$ws.Range("B3").Value2 = 0.104166666666667
$ws.Range("B3").NumberFormat = "HH:MM:SS"

# Move the active selection to the newly filled time cell, matching the
# author's final cursor position after entering the new row.
$ws.Range("B3").Select() | Out-Null
